# Apply scraped-price updates to the Adamantoise Profits workbook
# (generated from the authoritative diff; values set per-cell to match the
#  refreshed currentAveragePrice / profit calculations).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1131.6875
$ws.Range("I19").Value = 1173.625
$ws.Range("K19").Value = 1173.625
$ws.Range("M19").Value = -998.625

# Row 98
$ws.Range("H98").Value = 1543
$ws.Range("I98").Value = 1615.3572
$ws.Range("J98").Value = 1289.75
$ws.Range("K98").Value = 1615.3572
$ws.Range("L98").Value = 1289.75
$ws.Range("M98").Value = -117.3571999999999
$ws.Range("N98").Value = -4285.75

# Row 122
$ws.Range("H122").Value = 1543
$ws.Range("I122").Value = 1615.3572
$ws.Range("J122").Value = 1289.75
$ws.Range("K122").Value = 4846.071599999999
$ws.Range("L122").Value = 3869.25
$ws.Range("M122").Value = -2396.071599999999
$ws.Range("N122").Value = -8769.25

# Row 137
$ws.Range("H137").Value = 13891588
$ws.Range("I137").Value = 2999
$ws.Range("K137").Value = 8997
$ws.Range("M137").Value = -6447

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4349.684
$ws.Range("I61").Value = 4524.75
$ws.Range("J61").Value = 4222.364
$ws.Range("K61").Value = 4524.75
$ws.Range("L61").Value = 4222.364
$ws.Range("M61").Value = -4312.75
$ws.Range("N61").Value = -4646.364

# Row 74
$ws.Range("H74").Value = 2895.6667
$ws.Range("I74").Value = 3054.25
$ws.Range("K74").Value = 3054.25
$ws.Range("M74").Value = -2180.25

# Row 77
$ws.Range("H77").Value = 2895.6667
$ws.Range("I77").Value = 3054.25
$ws.Range("K77").Value = 15271.25
$ws.Range("M77").Value = -10903.25

# Row 136
$ws.Range("H136").Value = 4349.684
$ws.Range("I136").Value = 4524.75
$ws.Range("J136").Value = 4222.364
$ws.Range("K136").Value = 13574.25
$ws.Range("L136").Value = 12667.092
$ws.Range("M136").Value = -11024.25
$ws.Range("N136").Value = -17767.092

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3534.6943
$ws.Range("I86").Value = 3363.923
$ws.Range("K86").Value = 3363.923
$ws.Range("M86").Value = -2240.923

# Row 89
$ws.Range("H89").Value = 3534.6943
$ws.Range("I89").Value = 3363.923
$ws.Range("K89").Value = 16819.615
$ws.Range("M89").Value = -11203.615

$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Range("H33").Value = 7369.8335
$ws.Range("I33").Value = 7369.8335
$ws.Range("K33").Value = 7369.8335
$ws.Range("M33").Value = -6990.8335

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1099.6666
$ws.Range("J5").Value = 949.5
$ws.Range("L5").Value = 2848.5
$ws.Range("N5").Value = -3072.5

# Row 23
$ws.Range("H23").Value = 130.28572
$ws.Range("I23").Value = 330
$ws.Range("J23").Value = 97
$ws.Range("K23").Value = 990
$ws.Range("L23").Value = 291
$ws.Range("M23").Value = -755
$ws.Range("N23").Value = -761

# Row 107
$ws.Range("H107").Value = 666
$ws.Range("I107").Value = 668.5
$ws.Range("J107").Value = 664.2143
$ws.Range("K107").Value = 2005.5
$ws.Range("L107").Value = 1992.6429
$ws.Range("M107").Value = -85.5
$ws.Range("N107").Value = -5832.6429

# Row 120
$ws.Range("H120").Value = 17091.277
$ws.Range("I120").Value = 9182.556
$ws.Range("K120").Value = 27547.668
$ws.Range("M120").Value = -22709.668

# Row 135
$ws.Range("H135").Value = 1099.6666
$ws.Range("J135").Value = 949.5
$ws.Range("L135").Value = 8545.5
$ws.Range("N135").Value = -13615.5

$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 20000000
$ws.Range("I10").Value = 20000000
$ws.Range("K10").Value = 20000000
$ws.Range("M10").Value = -19999831

# Row 41
$ws.Range("H41").Value = 7065.25
$ws.Range("I41").Value = 6087.3335
$ws.Range("K41").Value = 6087.3335
$ws.Range("M41").Value = -5732.3335

# Row 48
$ws.Range("H48").Value = 15000
$ws.Range("J48").Value = 15000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15970

# Row 103
$ws.Range("H103").Value = 104624.625
$ws.Range("J103").Value = 104624.625
$ws.Range("L103").Value = 104624.625
$ws.Range("N103").Value = -106968.625

# Row 113
$ws.Range("H113").Value = 20986.54
$ws.Range("I113").Value = 15441.182
$ws.Range("K113").Value = 15441.182
$ws.Range("M113").Value = -13271.182

# Row 118
$ws.Range("H118").Value = 108989
$ws.Range("J118").Value = 108989
$ws.Range("L118").Value = 108989
$ws.Range("N118").Value = -112303

$ws = $wb.Worksheets.Item("LTW")
# Row 42
$ws.Range("H42").Value = 30039.5
$ws.Range("J42").Value = 30039.5
$ws.Range("L42").Value = 30039.5
$ws.Range("N42").Value = -31165.5

# Row 49
$ws.Range("H49").Value = 30039.5
$ws.Range("J49").Value = 30039.5
$ws.Range("L49").Value = 30039.5
$ws.Range("N49").Value = -30333.5

# Row 55
$ws.Range("H55").Value = 1022.6129
$ws.Range("I55").Value = 1032.619
$ws.Range("J55").Value = 1001.6
$ws.Range("K55").Value = 1032.619
$ws.Range("L55").Value = 1001.6
$ws.Range("M55").Value = -859.6189999999999
$ws.Range("N55").Value = -1347.6

# Row 61
$ws.Range("H61").Value = 4073.9
$ws.Range("I61").Value = 1392.3334
$ws.Range("J61").Value = 8096.25
$ws.Range("K61").Value = 1392.3334
$ws.Range("L61").Value = 8096.25
$ws.Range("M61").Value = -1190.3334
$ws.Range("N61").Value = -8500.25

# Row 113
$ws.Range("H113").Value = 4073.9
$ws.Range("I113").Value = 1392.3334
$ws.Range("J113").Value = 8096.25
$ws.Range("K113").Value = 1392.3334
$ws.Range("L113").Value = 8096.25
$ws.Range("M113").Value = 777.6666
$ws.Range("N113").Value = -12436.25

# Row 119
$ws.Range("H119").Value = 99894.664
$ws.Range("J119").Value = 99894.664
$ws.Range("L119").Value = 99894.664
$ws.Range("N119").Value = -109570.664

$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 57237.25
$ws.Range("J16").Value = 57237.25
$ws.Range("L16").Value = 57237.25
$ws.Range("N16").Value = -57821.25

# Row 47
$ws.Range("H47").Value = 55000
$ws.Range("J47").Value = 55000
$ws.Range("L47").Value = 55000
$ws.Range("N47").Value = -56144

# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# Row 132
$ws.Range("H132").Value = 3431.3462
$ws.Range("I132").Value = 3330.7856
$ws.Range("K132").Value = 9992.356800000001
$ws.Range("M132").Value = -7462.356800000001
